# Remove the "Chimères" group rows (rows 2-5) from the summary table.
# This shifts all subsequent rows ("Raies", "Requins") up by 4 rows,
# shrinking the used range from A1:D35 to A1:D31.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:D5").EntireRow.Delete() | Out-Null
